{"js": "// Adi\u00e7\u00e3o de escopo no projeto\n// Insert a new bulleted list item \"Emitir cupom fiscal ap\u00f3s pagamento\"\n// right after the existing \"Consulta de ve\u00edculos estacionados\" item.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst anchorText = \"Consulta de ve\u00edculos estacionados\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === anchorText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find anchor paragraph: \" + anchorText);\n}\n\n// Insert a new paragraph right after the anchor paragraph. This new\n// paragraph inherits the anchor's paragraph formatting (list style,\n// numbering, run formatting), matching the target diff.\nconst newPara = anchor.insertParagraph(\"Emitir cupom fiscal ap\u00f3s pagamento\", \"After\");\n\nawait context.sync();\n", "ps1": "# Adi\u00e7\u00e3o de escopo no projeto\n# Insert a new bulleted list item \"Emitir cupom fiscal ap\u00f3s pagamento\"\n# right after the existing \"Consulta de ve\u00edculos estacionados\" item.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"Consulta de ve\u00edculos estacionados\")\n\nif (-not $found) {\n    throw \"Could not find anchor paragraph: Consulta de ve\u00edculos estacionados\"\n}\n\n# Collapse to the end of the matched text, then insert a new paragraph\n# mark right after it. The new paragraph inherits the formatting\n# (ListParagraph style + bullet numbering) of the paragraph it was split\n# from, matching the target edit.\n$rng.Collapse(0)\n$rng.InsertParagraphAfter()\n\n# Move into the freshly created (still empty) paragraph and add the text.\n$rng.Collapse(0)\n$rng.Move(1, 1) | Out-Null\n$rng.InsertAfter(\"Emitir cupom fiscal ap\u00f3s pagamento\")\n"}
